# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2 of the
# per-language report sheets ("zh-cn" and "de-de").

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 19:06:08"
$wsZhCn.Range("H2").Value = "2016-03-21 19:06:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 19:06:12"
$wsDeDe.Range("H2").Value = "2016-03-21 19:06:53"
